$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (Price + Volume(1h) columns) with fresh market data.
# The Price column ("D") cells are stored as text (e.g. "534.30"), so force
# each modified cell back to Text format before writing the new value -- this
# stops Excel from auto-converting it to a number and dropping trailing zeros.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.156.51"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.519.55"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.30"
$ws.Range("E5").Value = "  -1.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.34"
$ws.Range("E6").Value = "  -3.87%  "
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.563"
$ws.Range("E8").Value = "  -2.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.523.98"
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.101"
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("E12").Value = "  -2.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.964.99"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.23"
$ws.Range("E15").Value = "  -2.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "59.087.51"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.518.56"
$ws.Range("E18").Value = "  -1.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.04"
$ws.Range("E19").Value = "  -2.06%  "
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.20"
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.82"
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.27"
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.426"
$ws.Range("E25").Value = "  -2.99%  "
$ws.Range("E27").Value = "  +0.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.80"
$ws.Range("E28").Value = "  -2.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.74"
$ws.Range("E29").Value = "  -2.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0773"
$ws.Range("E30").Value = "  -1.56%  "
$ws.Range("E31").Value = "  -2.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "164.74"
$ws.Range("E32").Value = "  +5.13%  "
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("E34").Value = "  -3.12%  "
$ws.Range("E35").Value = "  -7.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.49"
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("E37").Value = "  -3.58%  "
$ws.Range("E38").Value = "  -2.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.82"
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.67"
$ws.Range("E40").Value = "  -1.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.810"
$ws.Range("E41").Value = "  -2.44%  "
$ws.Range("E42").Value = "  -7.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "279.14"
$ws.Range("E43").Value = "  -6.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.598"
$ws.Range("E46").Value = "  -1.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0932"
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.41"
$ws.Range("E48").Value = "  -1.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.49"
$ws.Range("E49").Value = "  -1.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0512"
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("E51").Value = "  -2.02%  "
